$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 14:42"

# Rebuild the city data table (A4:E64), re-sorted descending by "Casos totales"
# and refreshed with updated case numbers.
$data = @(
    @('Madrid', 17166, 3882, 11194, 2090),
    @('Cataluña', 11592, 1697, 9223, 672),
    @('Castilla-La Mancha', 2780, 71, 2446, 263),
    @('Bizkaia/Vizcaya', 1850, 621, 1785, 65),
    @('Valencia/Valencia', 1788, 36, 1680, 72),
    @('Araba/Alava', 1435, 621, 1341, 94),
    @('Navarra', 1411, 35, 1327, 49),
    @('Ciudad Real', 1147, 95, 1050, 89),
    @('Alacant/Alicante', 1039, 17, 941, 81),
    @('La Rioja', 995, 48, 904, 43),
    @('Toledo', 965, 95, 860, 78),
    @('Aragon', 907, 29, 838, 40),
    @('Malaga', 905, 52, 811, 42),
    @('Zaragoza', 882, 42, 801, 39),
    @('A Coruña', 812, 47, 789, 23),
    @('Asturias', 779, 35, 719, 25),
    @('Pontevedra', 689, 47, 684, 5),
    @('Cantabria', 671, 14, 640, 17),
    @('Albacete', 666, 95, 592, 66),
    @('Gipuzkoa/Guipuzcoa', 661, 621, 640, 21),
    @('Granada', 648, 1, 615, 32),
    @('Caceres', 640, 24, 587, 35),
    @('Salamanca', 629, 73, 502, 54),
    @('Sevilla', 627, 8, 595, 24),
    @('Valladolid', 598, 50, 521, 27),
    @('Murcia', 596, 9, 579, 8),
    @('Tenerife', 539, 15, 519, 24),
    @('Burgos', 530, 67, 431, 32),
    @('Leon', 516, 48, 422, 46),
    @('Guadalajara', 428, 95, 370, 56),
    @('Castello/Castellon', 373, 4, 355, 14),
    @('Segovia', 361, 62, 262, 37),
    @('Cordoba', 349, 0, 342, 7),
    @('Jaen', 349, 5, 329, 15),
    @('Soria', 339, 32, 291, 16),
    @('Badajoz', 329, 24, 316, 4),
    @('Cadiz', 299, 4, 290, 5),
    @('Avila', 270, 53, 190, 27),
    @('Ourense', 235, 47, 231, 4),
    @('Mallorca', 210, 18, 194, 12),
    @('Cuenca', 177, 95, 142, 27),
    @('Gran Canaria', 171, 15, 169, 24),
    @('Palencia', 139, 14, 120, 5),
    @('Almeria', 134, 5, 122, 7),
    @('Lugo', 132, 47, 128, 4),
    @('Huesca', 119, 6, 109, 4),
    @('Teruel', 115, 6, 104, 5),
    @('Zamora', 106, 24, 73, 9),
    @('Huelva', 95, 2, 91, 2),
    @('Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
    @('Melilla', 39, 0, 38, 1),
    @('La Palma', 33, 15, 33, 24),
    @('Ibiza', 21, 18, 20, 1),
    @('Fuerteventura', 20, 15, 20, 24),
    @('Lanzarote', 17, 15, 17, 24),
    @('Ceuta', 16, 0, 16, 0),
    @('Menorca', 15, 18, 13, 0),
    @('Arroyo de la Luz', 7, 0, 7, 0),
    @('La Gomera', 4, 15, 2, 24),
    @('El Hierro', 3, 15, 3, 24),
    @('Formentera', 0, 10, 0, 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
    $ws.Cells.Item($r, 5).Value = $data[$i][4]
}
